# Adds a new "Method to use instead of deprecated executeMeRequestAsync"
# reference entry (with its Stack Overflow link) right before the trailing
# "_GoBack" bookmark at the end of the document, matching the existing
# "question / tab + link" paragraph-pair pattern used throughout the file.

$d = $word.ActiveDocument

# The trailing paragraph currently holds nothing but the (relocatable)
# "_GoBack" bookmark; it is the very last paragraph with real content
# above it, so grab it via the bookmark rather than a hard-coded index.
$bm = $d.Bookmarks.Item("_GoBack")
$target = $bm.Range.Paragraphs.Item(1)
$targetRange = $target.Range.Duplicate

# Flat-OPC "Word XML" fragment: replace that whole (currently-empty)
# paragraph with two paragraphs -- the first carries the new "Method to
# use instead of deprecated executeMeRequestAsync" text (split across two
# runs so a spell-check-skip proofErr can bracket the camelCase method
# name, mirroring the rest of the document), the second is a fresh
# tab + URL paragraph ending with the re-homed _GoBack bookmark, exactly
# like every other entry pair in the file.
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00D90E24" w:rsidRPr="0049162E" w:rsidRDefault="00D90E24" w:rsidP="0049162E"><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve">Method to use instead of deprecated </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>executeMeRequestAsync</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>http://stackoverflow.com/questions/18841084/how-to-undeprecate-facebook-code</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$targetRange.InsertXML($xml)
